# Performance & Air Quality Data
#
# Update a handful of Heart Rate (BPM) readings in column I, and leave the
# sheet scrolled/selected where the author was last working.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column I holds "Heart Rate (BPM)" readings.
$ws.Range("I13").Value = 157
$ws.Range("I22").Value = 160
$ws.Range("I33").Value = 149
$ws.Range("I37").Value = 159
$ws.Range("I41").Value = 160
$ws.Range("I45").Value = 162
$ws.Range("I46").Value = 151

# Scroll the viewport so row 16 is at the top, and leave I42 selected,
# matching where the author ended up when they saved the file.
$win = $excel.ActiveWindow
$win.ScrollRow = 16
$win.ScrollColumn = 1
$ws.Range("I42").Select()
